$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 1097
$ws.Range("I21").Value = 1097
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1097
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -629
$ws.Range("N21").Value = ""

$ws.Range("H23").Value = 1097
$ws.Range("I23").Value = 1097
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1097
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -863
$ws.Range("N23").Value = ""

$ws.Range("H64").Value = 3620.3076
$ws.Range("I64").Value = 3809.2856
$ws.Range("K64").Value = 3809.2856
$ws.Range("M64").Value = -3561.2856

$ws.Range("H67").Value = 3620.3076
$ws.Range("I67").Value = 3809.2856
$ws.Range("K67").Value = 3809.2856
$ws.Range("M67").Value = -2951.2856

$ws.Range("H70").Value = 13163.454
$ws.Range("I70").Value = 5716.5
$ws.Range("K70").Value = 17149.5
$ws.Range("M70").Value = -16879.5

$ws.Range("H73").Value = 13163.454
$ws.Range("I73").Value = 5716.5
$ws.Range("K73").Value = 17149.5
$ws.Range("M73").Value = -16213.5

$ws.Range("H75").Value = 80831.44500000001
$ws.Range("J75").Value = 99642.57000000001
$ws.Range("L75").Value = 99642.57000000001
$ws.Range("N75").Value = -101514.57

$ws.Range("H78").Value = 80831.44500000001
$ws.Range("J78").Value = 99642.57000000001
$ws.Range("L78").Value = 298927.71
$ws.Range("N78").Value = -308287.71

$ws.Range("H137").Value = 590382.8
$ws.Range("I137").Value = 910608
$ws.Range("J137").Value = 3303.3333
$ws.Range("K137").Value = 2731824
$ws.Range("L137").Value = 9909.999899999999
$ws.Range("M137").Value = -2729274
$ws.Range("N137").Value = -15009.9999

$ws.Range("H138").Value = 4313.6665
$ws.Range("I138").Value = 5207
$ws.Range("J138").Value = 3599
$ws.Range("K138").Value = 15621
$ws.Range("L138").Value = 10797
$ws.Range("M138").Value = -10481
$ws.Range("N138").Value = -21077

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 789
$ws.Range("I37").Value = 789
$ws.Range("K37").Value = 789
$ws.Range("M37").Value = -516

$ws.Range("H122").Value = 2028.1333
$ws.Range("I122").Value = 1965.7
$ws.Range("K122").Value = 5897.1
$ws.Range("M122").Value = -3447.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7279.8
$ws.Range("I86").Value = 653.1667
$ws.Range("K86").Value = 653.1667
$ws.Range("M86").Value = 469.8333

$ws.Range("H89").Value = 7279.8
$ws.Range("I89").Value = 653.1667
$ws.Range("K89").Value = 3265.8335
$ws.Range("M89").Value = 2350.1665

$ws.Range("H94").Value = 16279.6
$ws.Range("I94").Value = 7118
$ws.Range("K94").Value = 7118
$ws.Range("M94").Value = -6667

$ws.Range("H105").Value = 2747.6155
$ws.Range("I105").Value = 2235.9167
$ws.Range("K105").Value = 2235.9167
$ws.Range("M105").Value = -488.9167000000002

$ws.Range("H134").Value = 9448.789000000001
$ws.Range("I134").Value = 13767.556
$ws.Range("J134").Value = 5561.9
$ws.Range("K134").Value = 41302.66800000001
$ws.Range("L134").Value = 16685.7
$ws.Range("M134").Value = -38767.66800000001
$ws.Range("N134").Value = -21755.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3522.5
$ws.Range("I16").Value = 2490.4443
$ws.Range("K16").Value = 2490.4443
$ws.Range("M16").Value = -2203.4443

$ws.Range("H31").Value = 3693.25
$ws.Range("I31").Value = 2362.9
$ws.Range("K31").Value = 2362.9
$ws.Range("M31").Value = -2067.9

$ws.Range("H34").Value = 3693.25
$ws.Range("I34").Value = 2362.9
$ws.Range("K34").Value = 2362.9
$ws.Range("M34").Value = -2160.9

$ws.Range("H103").Value = 69999
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 69999
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 69999
$ws.Range("M103").Value = ""
$ws.Range("N103").Value = -72343

$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = ""

$ws.Range("H113").Value = 3522.5
$ws.Range("I113").Value = 2490.4443
$ws.Range("K113").Value = 2490.4443
$ws.Range("M113").Value = -320.4443000000001

$ws.Range("H134").Value = 5599.1377
$ws.Range("J134").Value = 8156
$ws.Range("L134").Value = 24468
$ws.Range("N134").Value = -29538

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 16966.5
$ws.Range("I86").Value = 600
$ws.Range("J86").Value = 33333
$ws.Range("K86").Value = 1800
$ws.Range("L86").Value = 99999
$ws.Range("M86").Value = -614
$ws.Range("N86").Value = -102371

$ws.Range("H89").Value = 16966.5
$ws.Range("I89").Value = 600
$ws.Range("J89").Value = 33333
$ws.Range("K89").Value = 5400
$ws.Range("L89").Value = 299997
$ws.Range("M89").Value = 528
$ws.Range("N89").Value = -311853

$ws.Range("H107").Value = 1630.7179
$ws.Range("J107").Value = 1661.1111
$ws.Range("L107").Value = 4983.3333
$ws.Range("N107").Value = -8823.3333

$ws.Range("H131").Value = 3243.48
$ws.Range("I131").Value = 1499.5
$ws.Range("J131").Value = 3316.1458
$ws.Range("K131").Value = 4498.5
$ws.Range("L131").Value = 9948.437399999999
$ws.Range("M131").Value = 541.5
$ws.Range("N131").Value = -20028.4374

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1432.3334
$ws.Range("I80").Value = 1432.3334
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 1432.3334
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -434.3334
$ws.Range("N80").Value = ""

$ws.Range("H83").Value = 1432.3334
$ws.Range("I83").Value = 1432.3334
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 7161.666999999999
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -2169.666999999999
$ws.Range("N83").Value = ""

$ws.Range("H97").Value = 1024.1111
$ws.Range("I97").Value = 235.28572
$ws.Range("J97").Value = 3785
$ws.Range("K97").Value = 235.28572
$ws.Range("L97").Value = 3785
$ws.Range("M97").Value = 260.71428
$ws.Range("N97").Value = -4777

$ws.Range("H102").Value = 5765
$ws.Range("I102").Value = 7022.143
$ws.Range("K102").Value = 7022.143
$ws.Range("M102").Value = -5400.143

$ws.Range("H122").Value = 1828.2222
$ws.Range("J122").Value = 2999.6667
$ws.Range("L122").Value = 8999.000100000001
$ws.Range("N122").Value = -13899.0001

$ws.Range("H132").Value = 4698.25
$ws.Range("I132").Value = 4740.8887
$ws.Range("K132").Value = 14222.6661
$ws.Range("M132").Value = -11692.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1741.8334
$ws.Range("I22").Value = 817
$ws.Range("K22").Value = 817
$ws.Range("M22").Value = -522

$ws.Range("H27").Value = 1741.8334
$ws.Range("I27").Value = 817
$ws.Range("K27").Value = 817
$ws.Range("M27").Value = -710

$ws.Range("H40").Value = 2095
$ws.Range("I40").Value = 1927
$ws.Range("K40").Value = 1927
$ws.Range("M40").Value = -1791

$ws.Range("H61").Value = 4997
$ws.Range("I61").Value = 4997
$ws.Range("K61").Value = 4997
$ws.Range("M61").Value = -4795

$ws.Range("H94").Value = 14165
$ws.Range("J94").Value = 14165
$ws.Range("L94").Value = 14165
$ws.Range("N94").Value = -15517

$ws.Range("H100").Value = 2582.5
$ws.Range("I100").Value = 2579
$ws.Range("K100").Value = 2579
$ws.Range("M100").Value = -2038

$ws.Range("H113").Value = 4997
$ws.Range("I113").Value = 4997
$ws.Range("K113").Value = 4997
$ws.Range("M113").Value = -2827

$ws.Range("H136").Value = 2219.8367
$ws.Range("I136").Value = 1302.0714
$ws.Range("K136").Value = 3906.2142
$ws.Range("M136").Value = -1356.2142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("L11").Value = 0
$ws.Range("N11").Value = ""

$ws.Range("H61").Value = 500
$ws.Range("I61").Value = 500
$ws.Range("K61").Value = 500
$ws.Range("M61").Value = -208

$ws.Range("H126").Value = 2723.2307
$ws.Range("I126").Value = 2700.1667
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 8100.500100000001
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -5630.500100000001
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 3498.5
$ws.Range("I132").Value = 2997
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 8991
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -6461
$ws.Range("N132").Value = -17060
